$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Arthur" no longer shares a cell with the device name "HP ENVY" (now on its own row)
$ws.Range("D4").Value = "Arthur"

# 2. Scenario 4.1 now asks the player to press 3 instead of 5
$ws.Range("B20").Value = "4.1 Appuyez sur 3 et ensuite Enter"

# 3. Column B (the long scenario-description column) got wider
$ws.Columns("B").ColumnWidth = 75.6666666667

# 4. Header row no longer needs the extra height for the two-line "Arthur" entry
$ws.Rows("4").EntireRow.AutoFit()

# 5. Move the active selection
$null = $ws.Range("G20").Select()
